# NumberFormat.xlsx - add two more "Number Format" example rows (fraction formats)
# and push the trailing "Some text" label down below them.
#
# Original layout:
#   Row 19: "Some text"  (text, style s=12, numFmtId 49 "@")
#
# New layout:
#   Row 19: 1.25  formatted as built-in fraction          "# ?/?"   (numFmtId 12)
#   Row 20: 1.25  formatted as custom fraction            "# ?/100"
#   Row 21: "Some text"  (moved down, same text style as before)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: turn the old "Some text" cell into a number using the built-in
# fraction number format ("# ?/?", numFmtId 12). Setting NumberFormat first
# (before the value) ensures the cell's old string type is cleared so the
# new numeric value is stored correctly.
$ws.Range("A19").NumberFormat = "# ?/?"
$ws.Range("A19").Value = 1.25

# Row 20 (new row): another fraction example, this time with a custom
# number format "# ?/100".
$ws.Range("A20").NumberFormat = "# ?/100"
$ws.Range("A20").Value = 1.25

# Row 21 (new row): re-add the "Some text" label that used to be on row 19,
# keeping its original general-text ("@") number format.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "Some text"

# Leave the selection on A19, matching where the edit was made.
[void]$ws.Range("A19").Select()
